# SurfaceWater.xlsx - "handle empty values for scsb model monthly values"
#
# 1. Update the 12 "Monthly discharge (m3/s)" template cells (row 19,
#    columns B:M) on the "Runoff and Streamflow Models" sheet so the
#    renderer pulls `.model_result`, falls back to a blank string via
#    `ifEmpty( )` when the model produced nothing, and only then rounds.
# 2. Make "Runoff and Streamflow Models" the active sheet/tab (it was
#    "Summary" before), with the selection sitting on C19, which also
#    nudges the sheet's used range down to row 20.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)   # "Runoff and Streamflow Models"

for ($i = 1; $i -le 12; $i++) {
    $col = $i + 1  # month 1 -> column B (2) ... month 12 -> column M (13)
    $new = "{ d.scsb2016_output.monthly_discharge.$i.model_result:ifEmpty( ):round(2) }"
    $ws3.Cells.Item(19, $col).Value = $new
}

# Grow the sheet's used range from A1:M19 to A1:M20 with a trailing blank
# row (same row height as the rest of the data rows), mirroring the extra
# empty <row r="20".../> left behind by the edit.
$ws3.Rows.Item(20).RowHeight = 15
$ws3.Cells.Item(20, 1).NumberFormat = "General"

# Switch the active sheet from "Summary" to "Runoff and Streamflow Models"
# and park the selection on C19.
$ws3.Range("C19").Select() | Out-Null
$ws3.Activate() | Out-Null
